$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 88.20920302487313
$ws.Range("C2").Value = 98.04764679952338
$ws.Range("D2").Value = 99.21846182535455
$ws.Range("E2").Value = 98.76026201257623
$ws.Range("F2").Value = 98.1892679370048
$ws.Range("G2").Value = 97.29729133519017
$ws.Range("H2").Value = 95.87286587321537

$ws.Range("B3").Value = 83.25910012305783
$ws.Range("C3").Value = 98.92298353786403
$ws.Range("D3").Value = 99.62945199591201
$ws.Range("E3").Value = 98.86846984517264
$ws.Range("F3").Value = 98.38721355232411
$ws.Range("G3").Value = 97.38471645842995
$ws.Range("H3").Value = 95.91265432175534

$ws.Range("B4").Value = 93.56253383708672
$ws.Range("C4").Value = 97.7253575610866
$ws.Range("D4").Value = 99.11420376632694
$ws.Range("E4").Value = 98.58637681750606
$ws.Range("F4").Value = 98.2571061269427
$ws.Range("G4").Value = 97.26772153246471
$ws.Range("H4").Value = 95.86440757321071

$ws.Range("B5").Value = 90.62446518683227
$ws.Range("C5").Value = 97.8019675006344
$ws.Range("D5").Value = 99.16963294961775
$ws.Range("E5").Value = 98.72429906188519
$ws.Range("F5").Value = 98.18816835193275
$ws.Range("G5").Value = 97.21114645633337
$ws.Range("H5").Value = 95.80088590653605

$ws.Range("B6").Value = 90.38861938304689
$ws.Range("C6").Value = 97.97436348302708
$ws.Range("D6").Value = 99.21650934838365
$ws.Range("E6").Value = 98.73276714353163
$ws.Range("F6").Value = 98.15703488507259
$ws.Range("G6").Value = 97.29854031452236
$ws.Range("H6").Value = 95.85249414651719
